# Update the per-country description cell (A1) on each scenario sheet with
# the shortened text, and restore the resulting autofit row height + cursor
# (selection) position on each sheet, finishing with "Vietnam" as the active
# sheet/tab (matching the workbook's stored activeTab).

$wb = $excel.ActiveWorkbook

$descriptions = @{
    "Australia" = "**2025 PV Production Costs in Australia across different scenarios**`n`nThis scenario illustrates the cost of final PV modules across 3 different scenarios:`n1) domestic manufacturing, `n2) imported wafers from China, `n3) imported cell from China.`n`nThe panel assembly is done domestically. The production capacity considered is 4 GW for wafer, cell and panel assembly."
    "Germany"   = "**2025 PV Production Costs in Germany across different scenarios**`n`nThis scenario illustrates the cost of final PV modules across 3 different scenarios:`n1) domestic manufacturing,`n2) imported wafers from China,`n3) imported cell from China.`n`nThe panel assembly is done domestically. The production capacity considered is 4 GW for wafer, cell and panel assembly.`n"
    "India"     = "**2025 PV Production Costs in India across different scenarios**`n`nThis scenario illustrates the cost of final PV modules across 3 different scenarios :`n1) domestic manufacturing, `n2) imported wafers from China,`n3) imported cell from China.`n`nThe panel assembly is done domestically. The production capacity considered is 4 GW for wafer, cell and panel assembly.`n"
    "Vietnam"   = "**2025 PV Production Costs in Vietnam across different scenarios**`n`nThis scenario illustrates the cost of final PV modules across 3 different scenarios :`n1) domestic manufacturing,`n2) imported wafers from China,`n3) imported cell from China. `n`nThe panel assembly is done domestically. The production capacity considered is 4 GW for wafer, cell and panel assembly.`n"
}

$rowHeights = @{
    "Australia" = 159
    "Germany"   = 174
    "India"     = 159
    "Vietnam"   = 174
}

$selections = @{
    "Australia" = "A4"
    "Germany"   = "A1"
    "India"     = "A4"
    "Vietnam"   = "C1"
}

# Order matters: visit the sheets in tab order and re-activate "Vietnam"
# last so it ends up the active tab/sheet again, same as before the edit.
$order = @("Australia", "Germany", "India", "Vietnam")

foreach ($name in $order) {
    $ws = $wb.Worksheets.Item($name)
    [void]$ws.Activate()

    $ws.Range("A1").Value = $descriptions[$name]
    $ws.Rows.Item(1).RowHeight = $rowHeights[$name]

    [void]$ws.Range($selections[$name]).Select()
}
